$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B (recognized text) and column C (count) for rows 2-14
$ws.Range("B2").Value = "<his>"
$ws.Range("C2").Value = 60

$ws.Range("B3").Value = "<when>"
$ws.Range("C3").Value = 61

$ws.Range("B4").Value = "<alt>"
$ws.Range("C4").Value = 60

$ws.Range("B5").Value = "<cal>"
$ws.Range("C5").Value = 58

$ws.Range("B6").Value = "<zo>"
$ws.Range("C6").Value = 62

$ws.Range("B7").Value = "<hes>"

$ws.Range("B8").Value = "<were>"

$ws.Range("B9").Value = "<somen>"
$ws.Range("C9").Value = 60

$ws.Range("B10").Value = "<come>"
$ws.Range("C10").Value = 62

$ws.Range("C11").Value = 62

$ws.Range("B12").Value = "<hese>"
$ws.Range("C12").Value = 62

$ws.Range("B13").Value = "<heshth>"

$ws.Range("B14").Value = "<ale>"
